# Update 200 runs descended data
# This script applies the data refresh that was captured for the
# "Optimizer | 200 run - Descended" worksheet (4th sheet in the workbook),
# and re-selects that sheet as the active tab/cell, matching the source
# commit "Update 200 runs descended data".

$wb = $excel.ActiveWorkbook

# Worksheets are in tab order:
# 1 Optimizer Disabled - Ascended
# 2 Optimizer Disabled - Descended
# 3 Optimizer | 200 run - Ascended
# 4 Optimizer | 200 run - Descended   <-- this is the sheet being updated
# 5 Optimizer | 1000 run - Ascended
# 6 Optimizer | 1000 run - Descend
# 7 Graphs
$ws = $wb.Worksheets.Item(4)

# New simulation results for column D (net income after refund) and
# column I (net income after refund for the other optimizer variant).
# Columns E, J, L, M, Q and R are formulas and recalculate automatically,
# as do the dependent chart caches on the Graphs sheet.
$rows = @(
    @{Row=3; D=320340; I=62266}
    @{Row=4; D=924402; I=179581}
    @{Row=5; D=1589312; I=338336}
    @{Row=6; D=2247622; I=526991}
    @{Row=7; D=2834458; I=722345}
    @{Row=8; D=2899332; I=745546}
    @{Row=9; D=2964140; I=769046}
    @{Row=10; D=3028882; I=792845}
    @{Row=11; D=3093558; I=816943}
    @{Row=12; D=3158168; I=841340}
    @{Row=13; D=3222712; I=866036}
    @{Row=14; D=3287190; I=891031}
    @{Row=15; D=3351602; I=916325}
    @{Row=16; D=3415948; I=941918}
    @{Row=17; D=3480228; I=967810}
    @{Row=18; D=3544442; I=994001}
    @{Row=19; D=4182952; I=1272356}
    @{Row=20; D=4814862; I=1580611}
    @{Row=21; D=5440172; I=1918766}
    @{Row=22; D=6058882; I=2286821}
    @{Row=23; D=6653668; I=2679976}
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 9).Value = $r.I
}

# Recalculate so the formula cells (E, J, L, M, Q, R) and the chart
# num-caches that read from this sheet pick up the new values.
$excel.Calculate()

# The author ended the edit with the "200 run - Descended" sheet active
# and cell F15 selected (previously "200 run - Ascended" / M17 had been
# the active tab/selection).
$ws.Activate()
$ws.Range("F15").Select() | Out-Null
